$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Vue.js, Spring Boot, MySQL"
# (a single run) without using Find/Replace (Replace triggers an unrelated
# whole-document re-serialization in this engine) - just use Find to
# locate the range, then edit it directly via the Range object.
$rng = $d.Content
$found = $rng.Find.Execute("Vue.js, Spring Boot, MySQL")
if (-not $found) {
    throw "Could not locate the stack-list paragraph text"
}

$start = $rng.Start
$end = $rng.End

# Rewrite the paragraph text into the new order: "MySQL, Spring Boot, Vue.js"
$rng.Text = "MySQL, Spring Boot, Vue.js"

# Boundaries within the combined text "MySQL, Spring Boot, Vue.js"
#   "MySQL" | ", Spring Boot, " | "Vue.js"
$b1 = $start + 5
$b2 = $start + 20
$end = $start + 26

# Force a run split at each boundary without leaving residual formatting
# behind, by toggling Bold off then back on (net value unchanged) across
# the trailing portion of the text - this causes the interop layer to
# materialize a new run boundary at $b1 / $b2.
$tail1 = $d.Range($b1, $end)
$tail1.Bold = 0
$tail1.Bold = 1

$tail2 = $d.Range($b2, $end)
$tail2.Bold = 0
$tail2.Bold = 1
